$wb = $excel.ActiveWorkbook

# samples_retained sheet: two corrupt files found -> decrement affected counts by 1 each
$ws1 = $wb.Worksheets.Item("samples_retained")
$ws1.Range("F2").Value = 485
$ws1.Range("G22").Value = 6429

# negative sheet: fix mislabeled emotion abbreviation (Contepmt row was "sad", should be "con")
$ws5 = $wb.Worksheets.Item("negative")
$ws5.Range("C18").Value = "con"

# leave selection where the user last clicked after editing F2
$ws1.Activate() | Out-Null
$ws1.Range("F3").Select() | Out-Null
